# Mise à jour de l'application
# Adds a new attendance column (CB) for the 2025-11-18 session, mirroring
# the existing CA (2025-11-15) column's layout/formatting, then updates
# the active-cell selection to reflect where the author left off editing.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. New header date in CB1 --------------------------------------------
$ws.Range("CB1").Value = 45979

# --- 2. New attendance marks for CB2:CB29 (row 12 has no CA/CB columns) --
$ws.Range("CB2").Value = "P"
$ws.Range("CB3").Value = "P"
$ws.Range("CB4").Value = "P"
$ws.Range("CB5").Value = "B"
$ws.Range("CB6").Value = "B"
$ws.Range("CB7").Value = "P"
$ws.Range("CB8").Value = "P"
$ws.Range("CB9").Value = "P"
$ws.Range("CB10").Value = "P"
$ws.Range("CB11").Value = "P"
$ws.Range("CB13").Value = "B"
$ws.Range("CB14").Value = "P"
$ws.Range("CB15").Value = "P"
$ws.Range("CB16").Value = "P"
$ws.Range("CB17").Value = "P"
$ws.Range("CB18").Value = "P"
$ws.Range("CB19").Value = "P"
$ws.Range("CB20").Value = "P"
$ws.Range("CB22").Value = "P"
$ws.Range("CB23").Value = "P"
$ws.Range("CB24").Value = "P"
$ws.Range("CB25").Value = "P"
$ws.Range("CB26").Value = "P"
$ws.Range("CB27").Value = "P"
$ws.Range("CB28").Value = "P"
$ws.Range("CB29").Value = "P"
# CB21 stays blank (that row has no data for any prior session either).

# --- 3. Copy formatting from column CA onto the new column CB ------------
# (done after the values are written so the per-cell recalculation that
# Value= triggers isn't short-circuited by the paste operation)
# Row 12's table stops at column AX (that player's row was never extended
# with CA/CB tracking cells), so it is deliberately skipped here.
for ($r = 1; $r -le 29; $r++) {
    if ($r -eq 12) { continue }
    $ws.Range("CA$r").Copy()
    $ws.Range("CB$r").PasteSpecial(-4122)
}

# --- 4. Restore the active selection left by the author ------------------
$ws.Range("CD28").Select()
